# Generate plans using AI
# Update rule values / labels used by the AI plan generator across the
# PaymentGoal, ConfidenceScore, Income and PaymentHistory rule sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# PaymentGoal sheet: tweak Tenure adjustments
# ---------------------------------------------------------------------
$wsPaymentGoal = $wb.Worksheets.Item("PaymentGoal")
$wsPaymentGoal.Activate() | Out-Null
$wsPaymentGoal.Range("B2").Value = 0
$wsPaymentGoal.Range("B3").Value = -6
$wsPaymentGoal.Range("B6").Select() | Out-Null

# ---------------------------------------------------------------------
# Income sheet: rename header "Tenure" -> "Rate" and tweak values
# ---------------------------------------------------------------------
$wsIncome = $wb.Worksheets.Item("Income")
$wsIncome.Activate() | Out-Null
$wsIncome.Range("B1").Value = "Rate"
$wsIncome.Range("B2").Value = 0.7
$wsIncome.Range("B4").Value = -0.7
$wsIncome.Range("G9").Select() | Out-Null

# ---------------------------------------------------------------------
# PaymentHistory sheet: rename "MOSTLY DISCIPLINED" -> "MOSTLY_DISCIPLINED"
# and make it the active sheet/selection
# ---------------------------------------------------------------------
$wsPaymentHistory = $wb.Worksheets.Item("PaymentHistory")
$wsPaymentHistory.Activate() | Out-Null
$wsPaymentHistory.Range("A4").Value = "MOSTLY_DISCIPLINED"
$wsPaymentHistory.Range("A4").Select() | Out-Null
